$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value2 = 1
$ws.Cells.Item($row, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value2 = "Arica y Parinacota"

$ws.Cells.Item($row, 4).Value2 = 44890
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value2 = 15
$ws.Cells.Item($row, 6).Value2 = 100114007
$ws.Cells.Item($row, 7).Value2 = "Jengibre"
$ws.Cells.Item($row, 8).Value2 = "Sin especificar"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 300
$ws.Cells.Item($row, 11).Value2 = 14000
$ws.Cells.Item($row, 12).Value2 = 15000
$ws.Cells.Item($row, 13).Value2 = 14500
$ws.Cells.Item($row, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item($row, 15).Value2 = "Perú"
$ws.Cells.Item($row, 16).Value2 = 1115
$ws.Cells.Item($row, 17).Value2 = 13
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
